$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 54, pushing existing rows 54-74 down to 55-75.
$ws.Rows.Item(54).Insert()

# Populate the new row 54 with the new weekly record.
$ws.Range("A54").Value2 = 11
$ws.Range("B54").Value2 = "Vega Monumental Concepción"
$ws.Range("C54").Value2 = "Bíobío"
$ws.Range("D54").Value2 = 44518
$ws.Range("E54").Value2 = 8
$ws.Range("F54").Value2 = 100112032
$ws.Range("G54").Value2 = "Zapallo italiano"
$ws.Range("H54").Value2 = "Sin especificar"
$ws.Range("I54").Value2 = "Primera"
$ws.Range("J54").Value2 = 450
$ws.Range("K54").Value2 = 6500
$ws.Range("L54").Value2 = 7000
$ws.Range("M54").Value2 = 6722
$ws.Range("N54").Value2 = "$/caja 60 unidades"
$ws.Range("O54").Value2 = "Región de O'Higgins"
$ws.Range("P54").Value2 = 112
$ws.Range("Q54").Value2 = 60
$ws.Range("R54").Value2 = "Hortaliza"
